# Update column F ("dSF") values per repulled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 2
    4  = 3
    5  = -1
    8  = 3
    9  = -1
    10 = 3
    11 = -1
    12 = 2
    13 = -3
    14 = -1
    16 = -5
    17 = -2
    19 = 9
    20 = 3
    21 = 1
    22 = 6
    23 = -1
    24 = 2
    25 = 2
    26 = -4
    27 = 12
    28 = 0
    29 = -1
    30 = 1
    31 = 4
    32 = 1
    33 = 3
    34 = 8
    35 = 1
    37 = -1
    38 = -1
    39 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
